$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.389.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.627.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "325.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("E7").Value = "  -0.97%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0810"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.20%  "

$ws.Range("E13").Value = "  +1.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.050.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.641.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.851"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.418.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "267.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.33%  "

$ws.Range("E25").Value = "  -1.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "25.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.137"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0806"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.84%  "

$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.60%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.111"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.09%  "

$ws.Range("E44").Value = "  +5.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.057.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.30%  "

$ws.Range("E48").Value = "  -5.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.86%  "
